# Auto-generated edit script: Add data for 2024-02-17
# Updates 2024 (column K, and a couple spillover column J) crime counts
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 867
$ws.Range("J3").Value = 8074
$ws.Range("K3").Value = 816
$ws.Range("K4").Value = 187
$ws.Range("K5").Value = 48
$ws.Range("K6").Value = 1137
$ws.Range("J7").Value = 29253
$ws.Range("K7").Value = 3055

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 22
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 91
$ws.Range("K8").Value = 184
$ws.Range("K9").Value = 16
$ws.Range("K10").Value = 19
$ws.Range("K11").Value = 66
$ws.Range("K15").Value = 19
$ws.Range("K18").Value = 22
$ws.Range("K19").Value = 79
$ws.Range("K20").Value = 70
$ws.Range("K26").Value = 6
$ws.Range("K29").Value = 151
$ws.Range("K30").Value = 5
$ws.Range("K33").Value = 130
$ws.Range("K34").Value = 21
$ws.Range("K36").Value = 36
$ws.Range("K37").Value = 93
$ws.Range("K40").Value = 5
$ws.Range("K41").Value = 30
$ws.Range("K42").Value = 98
$ws.Range("K43").Value = 31
$ws.Range("K47").Value = 20
$ws.Range("K51").Value = 44
$ws.Range("K52").Value = 78
$ws.Range("K53").Value = 38
$ws.Range("K54").Value = 55
$ws.Range("K55").Value = 30
$ws.Range("K63").Value = 13
$ws.Range("K67").Value = 132
$ws.Range("K71").Value = 9
$ws.Range("J72").Value = 110
$ws.Range("K73").Value = 34
$ws.Range("K76").Value = 45
$ws.Range("K77").Value = 18
$ws.Range("K78").Value = 40
$ws.Range("K79").Value = 82
$ws.Range("K83").Value = 57
$ws.Range("K84").Value = 26
$ws.Range("K85").Value = 152
$ws.Range("K86").Value = 23
$ws.Range("K87").Value = 2
$ws.Range("K88").Value = 40
$ws.Range("K89").Value = 50
$ws.Range("K95").Value = 52
$ws.Range("K96").Value = 47
$ws.Range("K98").Value = 18
$ws.Range("K99").Value = 60
$ws.Range("J101").Value = 29253
$ws.Range("K101").Value = 3055

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 39
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 38
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K3").Value = 7
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 9
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 38
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K2").Value = 4
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 11
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 10
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 47
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 19
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 2
